# Fruta / hortaliza, semanal
# Insert a new weekly record at row 27 (pushing existing rows 27-103 down to 28-104)
# and populate it with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 27; this shifts rows 27-103 down to 28-104
# and extends the sheet dimension from A1:R103 to A1:R104.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new weekly record.
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = 'Vega Monumental Concepción'
$ws.Range("C27").Value = 'Bíobío'
$ws.Range("D27").Value = 44624
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112032
$ws.Range("G27").Value = 'Zapallo italiano'
$ws.Range("H27").Value = 'Sin especificar'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 11000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 11467
$ws.Range("N27").Value = '$/caja 60 unidades'
$ws.Range("O27").Value = 'Región Metropolitana'
$ws.Range("P27").Value = 191
$ws.Range("Q27").Value = 60
$ws.Range("R27").Value = 'Hortaliza'

# Make sure the date cell keeps the same date number format used by the rest of column D.
$ws.Range("D27").NumberFormat = $ws.Range("D28").NumberFormat
